$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 21:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1441742
$ws.Range("C4").Value = 11394
$ws.Range("E4").Value = 1044506
$ws.Range("G4").Value = 318
$ws.Range("H4").Value = 85515

# Row 10 - Francia
$ws.Range("B10").Value = 178870
$ws.Range("C10").Value = 810
$ws.Range("E10").Value = 91840

# Row 11 - Alemania
$ws.Range("B11").Value = 174584
$ws.Range("C11").Value = 486
$ws.Range("E11").Value = 16395
$ws.Range("F11").Value = 1329
$ws.Range("G11").Value = 28
$ws.Range("H11").Value = 7889

# Row 15 - India
$ws.Range("D15").Value = 27969
$ws.Range("E15").Value = 51372

# Row 114 - Principado de Andorra
$ws.Range("B114").Value = 761
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 596
$ws.Range("E114").Value = 116
